$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows (2..14) down to (3..15)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new conference entry
$ws.Range("A2").Value = "针对中国年轻人控烟广告有潜力的信念"
$ws.Range("B2").Value = "裴瑞，于莲，赵亮，陈静茜"
$ws.Range("C2").Value = 2020
$ws.Range("D2").Value = "November"
$ws.Range("E2").Value = "The Medicine, Humanity and Media: Health China & Health Communication."
$ws.Range("F2").Value = "Online/Beijing"
$ws.Range("H2").Value = "https://www.bilibili.com/video/BV1Uv4116737"

# Match the final active-cell selection state recorded in the saved workbook
$ws.Range("H12").Select() | Out-Null
